$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "329.17" }
    @{ Cell = "E2"; Value = "1.34%" }
    @{ Cell = "D3"; Value = "41.17" }
    @{ Cell = "E3"; Value = "4.15%" }
    @{ Cell = "D4"; Value = "5.618" }
    @{ Cell = "E4"; Value = "-1.29%" }
    @{ Cell = "D5"; Value = "0.08196" }
    @{ Cell = "E5"; Value = "2.41%" }
    @{ Cell = "D6"; Value = "8.753" }
    @{ Cell = "E6"; Value = "1.50%" }
    @{ Cell = "D7"; Value = "2.001" }
    @{ Cell = "E7"; Value = "-0.29%" }
    @{ Cell = "D8"; Value = "4.492" }
    @{ Cell = "E8"; Value = "-0.04%" }
    @{ Cell = "E9"; Value = "-0.93%" }
    @{ Cell = "D10"; Value = "0.9215" }
    @{ Cell = "E10"; Value = "-0.41%" }
    @{ Cell = "D11"; Value = "0.1282" }
    @{ Cell = "E11"; Value = "3.37%" }
    @{ Cell = "D12"; Value = "0.1954" }
    @{ Cell = "E12"; Value = "-1.23%" }
    @{ Cell = "D13"; Value = "0.09365" }
    @{ Cell = "E13"; Value = "1.09%" }
    @{ Cell = "D14"; Value = "0.03891" }
    @{ Cell = "E14"; Value = "7.23%" }
    @{ Cell = "E15"; Value = "0.90%" }
    @{ Cell = "D16"; Value = "0.001298" }
    @{ Cell = "E16"; Value = "0.56%" }
    @{ Cell = "D17"; Value = "0.006105" }
    @{ Cell = "E17"; Value = "0.33%" }
    @{ Cell = "D19"; Value = "3.446" }
    @{ Cell = "E19"; Value = "2.83%" }
    @{ Cell = "E20"; Value = "0.28%" }
    @{ Cell = "D21"; Value = "8.246" }
    @{ Cell = "E21"; Value = "-5.45%" }
    @{ Cell = "E23"; Value = "0.03%" }
    @{ Cell = "D24"; Value = "0.04397" }
    @{ Cell = "E24"; Value = "-0.23%" }
    @{ Cell = "D25"; Value = "0.001258" }
    @{ Cell = "E25"; Value = "-0.15%" }
    @{ Cell = "D26"; Value = "0.004310" }
    @{ Cell = "E26"; Value = "-7.78%" }
    @{ Cell = "E27"; Value = "4.32%" }
    @{ Cell = "D39"; Value = "0.02787" }
    @{ Cell = "E39"; Value = "11.66%" }
    @{ Cell = "D40"; Value = "0.05400" }
    @{ Cell = "E40"; Value = "1.38%" }
    @{ Cell = "D41"; Value = "0.007807" }
    @{ Cell = "E41"; Value = "3.92%" }
    @{ Cell = "D42"; Value = "0.1417" }
    @{ Cell = "D43"; Value = "0.008943" }
    @{ Cell = "E43"; Value = "-6.74%" }
    @{ Cell = "D44"; Value = "0.002171" }
    @{ Cell = "E44"; Value = "2.56%" }
    @{ Cell = "D45"; Value = "0.01158" }
    @{ Cell = "E45"; Value = "1.34%" }
    @{ Cell = "D46"; Value = "0.00006760" }
    @{ Cell = "E46"; Value = "0.66%" }
    @{ Cell = "D47"; Value = "0.00000000750" }
    @{ Cell = "E47"; Value = "0.05%" }
    @{ Cell = "D48"; Value = "0.003192" }
    @{ Cell = "E48"; Value = "7.47%" }
    @{ Cell = "D49"; Value = "0.002281" }
    @{ Cell = "E49"; Value = "-0.42%" }
    @{ Cell = "D50"; Value = "0.00002101" }
    @{ Cell = "E50"; Value = "0.05%" }
    @{ Cell = "D51"; Value = "0.0002001" }
    @{ Cell = "E51"; Value = "0.05%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
